$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.773.28"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "2.102.84"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "227.81"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.31"
$ws.Range("E7").Value = "  +1.56%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +2.06%  "

$ws.Range("E10").Value = "  -0.39%  "

$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("D12").Value = "15.75"
$ws.Range("E12").Value = "  +6.35%  "

$ws.Range("D13").Value = "2.414.65"
$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("E14").Value = "  -1.31%  "

$ws.Range("E15").Value = "  +3.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.53"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("D17").Value = "2.105.64"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").Value = "38.761.18"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.90"
$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.76"
$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  -3.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("E26").Value = "  +2.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.33"
$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("E28").Value = "  +4.53%  "

$ws.Range("E29").Value = "  +5.06%  "

$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("D31").Value = "2.52"
$ws.Range("E31").Value = "  +10.24%  "

$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("E34").Value = "  -0.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.01"
$ws.Range("E35").Value = "  +7.34%  "

$ws.Range("E36").Value = "  +2.05%  "

$ws.Range("E37").Value = "  +0.75%  "

$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.14"
$ws.Range("E40").Value = "  -2.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.74"
$ws.Range("E41").Value = "  +2.89%  "

$ws.Range("E42").Value = "  +4.03%  "

$ws.Range("D43").Value = "1.526.11"
$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("E44").Value = "  +7.10%  "

$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("D46").Value = "7.79"
$ws.Range("E46").Value = "  +1.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0912"
$ws.Range("E47").Value = "  -0.26%  "

$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("E49").Value = "  +1.93%  "

$ws.Range("E50").Value = "  -0.62%  "

$ws.Range("D51").Value = "2.300.29"
$ws.Range("E51").Value = "  +0.28%  "
